function Set-TextValue {
    param($Ws, $Addr, $Val)
    $cell = $Ws.Range($Addr)
    $cell.NumberFormat = "@"
    $cell.Value2 = $Val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '27.935.36'
Set-TextValue $ws 'E2' '  +0.54%  '
Set-TextValue $ws 'D3' '1.812.44'
Set-TextValue $ws 'E4' '  -0.03%  '
Set-TextValue $ws 'D5' '310.71'
Set-TextValue $ws 'E5' '  +0.07%  '
Set-TextValue $ws 'E6' '  -0.02%  '
Set-TextValue $ws 'D7' '0.4973'
Set-TextValue $ws 'E7' '  -3.32%  '
Set-TextValue $ws 'D8' '0.3936'
Set-TextValue $ws 'E8' '  +3.86%  '
Set-TextValue $ws 'D9' '0.09670'
Set-TextValue $ws 'E9' '  +24.59%  '
Set-TextValue $ws 'D10' '1.103'
Set-TextValue $ws 'E10' '  +1.80%  '
Set-TextValue $ws 'D11' '40.97'
Set-TextValue $ws 'E11' '  -0.44%  '
Set-TextValue $ws 'D12' '6.450'
Set-TextValue $ws 'E12' '  +4.09%  '
Set-TextValue $ws 'D13' '20.48'
Set-TextValue $ws 'E13' '  +1.89%  '
Set-TextValue $ws 'E14' '  -0.06%  '
Set-TextValue $ws 'D15' '1.816.43'
Set-TextValue $ws 'E15' '  +2.22%  '
Set-TextValue $ws 'D16' '7.276'
Set-TextValue $ws 'E16' '  +1.69%  '
Set-TextValue $ws 'D17' '0.00001128'
Set-TextValue $ws 'E17' '  +5.35%  '
Set-TextValue $ws 'D18' '92.46'
Set-TextValue $ws 'E18' '  +1.13%  '
Set-TextValue $ws 'D19' '0.06665'
Set-TextValue $ws 'E19' '  +1.85%  '
Set-TextValue $ws 'D20' '1.002'
Set-TextValue $ws 'E20' '  +0.00%  '
Set-TextValue $ws 'D21' '17.13'
Set-TextValue $ws 'E21' '  +0.95%  '
Set-TextValue $ws 'D22' '5.919'
Set-TextValue $ws 'E22' '  +0.28%  '
Set-TextValue $ws 'D23' '27.986.80'
Set-TextValue $ws 'E23' '  +0.53%  '
Set-TextValue $ws 'D24' '11.15'
Set-TextValue $ws 'E24' '  +1.47%  '
Set-TextValue $ws 'D25' '2.246'
Set-TextValue $ws 'E25' '  +0.24%  '
Set-TextValue $ws 'D26' '159.66'
Set-TextValue $ws 'E26' '  +0.48%  '
Set-TextValue $ws 'D27' '2.023.82'
Set-TextValue $ws 'E27' '  +2.03%  '
Set-TextValue $ws 'D28' '20.57'
Set-TextValue $ws 'E28' '  +1.73%  '
Set-TextValue $ws 'E29' '  +1.40%  '
Set-TextValue $ws 'D30' '128.13'
Set-TextValue $ws 'E30' '  +2.31%  '
Set-TextValue $ws 'D31' '0.1064'
Set-TextValue $ws 'E31' '  -0.87%  '
Set-TextValue $ws 'D32' '1.036'
Set-TextValue $ws 'E32' '  +0.82%  '
Set-TextValue $ws 'D33' '5.573'
Set-TextValue $ws 'E33' '  +1.77%  '
Set-TextValue $ws 'D34' '3.640'
Set-TextValue $ws 'E34' '  +0.61%  '
Set-TextValue $ws 'D35' '0.06715'
Set-TextValue $ws 'E35' '  -5.41%  '
Set-TextValue $ws 'D36' '8.952'
Set-TextValue $ws 'E36' '  +3.14%  '
Set-TextValue $ws 'D37' '0.02322'
Set-TextValue $ws 'E37' '  +0.40%  '
Set-TextValue $ws 'D38' '0.2139'
Set-TextValue $ws 'E38' '  +0.85%  '
Set-TextValue $ws 'D39' '4.942'
Set-TextValue $ws 'E39' '  -1.31%  '
Set-TextValue $ws 'D40' '11.24'
Set-TextValue $ws 'E40' '  -2.71%  '
Set-TextValue $ws 'D41' '0.6178'
Set-TextValue $ws 'E41' '  +1.70%  '
Set-TextValue $ws 'D42' '1.002'
Set-TextValue $ws 'E42' '  +0.00%  '
Set-TextValue $ws 'D43' '1.148'
Set-TextValue $ws 'E43' '  -0.24%  '
Set-TextValue $ws 'D44' '13.14'
Set-TextValue $ws 'E44' '  +0.61%  '
Set-TextValue $ws 'D45' '0.5886'
Set-TextValue $ws 'E45' '  -0.83%  '
Set-TextValue $ws 'B46' 'WEMIXTOKEN'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws 'D46' '1.287'
Set-TextValue $ws 'E46' '  -2.56%  '
Set-TextValue $ws 'B47' 'PancakeSwap'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws 'D47' '3.692'
Set-TextValue $ws 'E47' '  -0.48%  '
Set-TextValue $ws 'D48' '123.07'
Set-TextValue $ws 'E48' '  -3.69%  '
Set-TextValue $ws 'E49' '  +2.27%  '
Set-TextValue $ws 'E50' '  -3.00%  '
Set-TextValue $ws 'D51' '0.06784'
Set-TextValue $ws 'E51' '  +0.78%  '
